# This script inserts a new weekly price record at row 15 of the sheet
# (pushing the existing rows 15-69 down to 16-70), matching the commit
# "Fruta / hortaliza, semanal" which adds a new data row for
# Agrícola del Norte S.A. de Arica - Papa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15; this shifts rows 15-69 down to 16-70 and keeps
# their data/formatting intact.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data.
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = Get-Date -Year 2022 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100114001
$ws.Cells.Item(15, 7).Value = "Papa"
$ws.Cells.Item(15, 8).Value = "Asterix"
$ws.Cells.Item(15, 9).Value = "1a (cosecha)"
$ws.Cells.Item(15, 10).Value = 1000
$ws.Cells.Item(15, 11).Value = 8500
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 8750
$ws.Cells.Item(15, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(15, 16).Value = 350
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
